# "upgrade left table until javakheti" - add a 2023 (column K) data point to
# the left-hand table on the Tianeti sheet, matching the formatting already
# used for the 2022 column (J).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone column J's formatting (number format / borders / alignment) onto
# column K so the new cells render exactly like the rest of the table.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats

# Header year + the three data rows (total / women / men).
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 404.8
$ws.Range("K5").Value = 268
$ws.Range("K6").Value = 534.9
